# Apply the cryptos-list refresh described by the commit
# ("Updated cryptos list on Fri Mar 31 19:02:33 UTC 2023 with GitHub Actions").
# Row 2-51 hold one crypto-coin each (B=Coin, C=Link, D=Price, E=Volume(1h));
# this pass rewrites the Price/Volume columns with the latest scrape, and for
# two pairs of rows (17/18 and 38/39) the coins also swapped rank/position.
#
# Numeric-looking price strings (e.g. "1.001", "41.82") are written with a
# leading apostrophe so Excel stores them as text (matching the original
# inline-string cells) instead of silently coercing them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.341.06'
$ws.Range("E2").Value = '  +1.53%  '

$ws.Range("D3").Value = '1.824.63'
$ws.Range("E3").Value = '  +2.66%  '

$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '''317.32'
$ws.Range("E5").Value = '  +0.63%  '

$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").Value = '''0.5334'
$ws.Range("E7").Value = '  -1.00%  '

$ws.Range("D8").Value = '''0.4033'
$ws.Range("E8").Value = '  +7.28%  '

$ws.Range("D9").Value = '''0.07592'
$ws.Range("E9").Value = '  +2.08%  '

$ws.Range("D10").Value = '''41.82'
$ws.Range("E10").Value = '  +0.48%  '

$ws.Range("D11").Value = '''1.106'
$ws.Range("E11").Value = '  +1.37%  '

$ws.Range("D12").Value = '''6.314'
$ws.Range("E12").Value = '  +4.11%  '

$ws.Range("D13").Value = '''1.001'
$ws.Range("E13").Value = '  +0.00%  '

$ws.Range("D14").Value = '''7.606'
$ws.Range("E14").Value = '  +5.72%  '

$ws.Range("D15").Value = '''20.78'
$ws.Range("E15").Value = '  +1.79%  '

$ws.Range("D16").Value = '1.826.13'
$ws.Range("E16").Value = '  +3.38%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '''0.00001073'
$ws.Range("E17").Value = '  +1.96%  '

$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Value = '''89.26'
$ws.Range("E18").Value = '  +1.46%  '

$ws.Range("D19").Value = '''0.06606'
$ws.Range("E19").Value = '  +2.88%  '

$ws.Range("D20").Value = '''17.65'
$ws.Range("E20").Value = '  +2.53%  '

$ws.Range("E21").Value = '  -0.01%  '

$ws.Range("D22").Value = '''6.102'
$ws.Range("E22").Value = '  +3.91%  '

$ws.Range("D23").Value = '28.371.76'
$ws.Range("E23").Value = '  +1.51%  '

$ws.Range("D24").Value = '''11.16'
$ws.Range("E24").Value = '  +0.33%  '

$ws.Range("D25").Value = '''2.200'
$ws.Range("E25").Value = '  +5.74%  '

$ws.Range("D26").Value = '''2.454'
$ws.Range("E26").Value = '  +8.07%  '

$ws.Range("D27").Value = '''158.20'
$ws.Range("E27").Value = '  +1.52%  '

$ws.Range("D28").Value = '''20.56'
$ws.Range("E28").Value = '  +1.70%  '

$ws.Range("D29").Value = '2.036.36'
$ws.Range("E29").Value = '  +3.10%  '

$ws.Range("D30").Value = '''123.76'
$ws.Range("E30").Value = '  +3.42%  '

$ws.Range("D31").Value = '''1.121'
$ws.Range("E31").Value = '  +1.05%  '

$ws.Range("D32").Value = '''0.1100'
$ws.Range("E32").Value = '  +4.30%  '

$ws.Range("D33").Value = '''5.652'
$ws.Range("E33").Value = '  +2.70%  '

$ws.Range("D34").Value = '''3.646'
$ws.Range("E34").Value = '  +0.12%  '

$ws.Range("D35").Value = '''0.07365'
$ws.Range("E35").Value = '  +15.69%  '

$ws.Range("D36").Value = '''0.2231'
$ws.Range("E36").Value = '  -0.67%  '

$ws.Range("D37").Value = '''0.02338'
$ws.Range("E37").Value = '  +3.44%  '

$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = '''5.191'
$ws.Range("E38").Value = '  +4.74%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '''8.833'
$ws.Range("E39").Value = '  +5.33%  '

$ws.Range("D40").Value = '''0.6250'
$ws.Range("E40").Value = '  +2.15%  '

$ws.Range("D41").Value = '''11.27'
$ws.Range("E41").Value = '  +2.46%  '

$ws.Range("E42").Value = '  +0.50%  '

$ws.Range("E43").Value = '  +0.01%  '

$ws.Range("E44").Value = '  -1.99%  '

$ws.Range("E45").Value = '  +2.06%  '

$ws.Range("E46").Value = '  +1.36%  '

$ws.Range("D47").Value = '''0.5835'
$ws.Range("E47").Value = '  +1.78%  '

$ws.Range("D48").Value = '''125.06'
$ws.Range("E48").Value = '  -0.63%  '

$ws.Range("D49").Value = '''1.989'
$ws.Range("E49").Value = '  +3.63%  '

$ws.Range("D50").Value = '''1.202'
$ws.Range("E50").Value = '  +1.54%  '

$ws.Range("D51").Value = '''0.06892'
$ws.Range("E51").Value = '  +1.62%  '
